# Adds a "Received" column to Table1 (the PCB/order tracking table in the
# top block of Sheet1), positioned just before the existing "Notes" column,
# and records that the order in row 6 (Petal v0.0 hardware order) was
# received on Friday, Octoder 18, 2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)   # Table1: B3:H11

# --- Move the existing "Notes" column (H) one slot to the right (I) -------
# Copy preserves both values and formatting (number format / style), then
# we clear the old column's data rows (header will be overwritten below).
$ws.Range("H3:H11").Copy($ws.Range("I3:I11"))
$ws.Range("H4:H11").ClearContents()

# --- Grow the table to cover the new column --------------------------------
$lo.Resize($ws.Range("B3:I11"))

# --- Fix up the header text for both columns (must happen after Resize so
# the table metadata picks up the rename) -----------------------------------
$ws.Range("H3").Value = "Received"
$ws.Range("I3").Value = "Notes"

# --- Populate the new "Received" column ------------------------------------
# Row 6 = order #3 "Petal v0.0 hardware order" -> received Fri Oct 18 2024.
$ws.Range("H6").Value = "Friday,Octoder 18, 2024"
$ws.Range("H6").NumberFormat = "[$-F800]dddd, mmmm dd, yyyy"

# Match the row height bump that came from Excel autofitting the new text.
$ws.Rows.Item(6).RowHeight = 31

# The rest of the new "Received" data cells (rows 4,5,7-11) stay blank; row 5
# picks up the short-date number format used elsewhere in the sheet for
# blank date cells.
$ws.Range("H5").NumberFormat = "m/d/yyyy"

Write-Host "Received column added to Table1"
